# Apply the commit: split the greeting sentence into proofed runs
# ("Anzhelika" / "Rusyatinskaya" wrapped in spellcheck proofErr markers)
# and append a new paragraph "I created this file yesterday." after it,
# moving the _GoBack bookmark to the end of the new paragraph - the way
# Word itself would leave things after you type the extra sentence.

$d = $word.ActiveDocument

# --- Step 1: grow the text in place (no paragraph break yet) ---------
# Doing this via Find/Replace (rather than touching Range boundaries
# directly) keeps the existing (collapsed) _GoBack bookmark glued to
# the very end of the story, where Word leaves it after an edit.
$find = $d.Content.Find
$find.Execute(
    "Rusyatinskaya.", $true, $false, $false, $false, $false, $true, 1,
    $false, "Rusyatinskaya.I created this file yesterday.", 2)

# --- Step 2: split that single paragraph into two --------------------
# The break goes right after "...Rusyatinskaya." (position 35 - the
# length of "My name is Anzhelika Rusyatinskaya."). Splitting this way
# (rather than re-building the whole range) lets _GoBack naturally ride
# along to the end of the newly typed second paragraph.
$splitPoint = $d.Range(35, 35)
$splitPoint.InsertParagraphAfter()

# --- Step 3: mark the name as spell-checked in the first paragraph ---
# Rebuild just paragraph 1's text (now isolated from the bookmark,
# which followed paragraph 2) into separate runs bracketed by
# <w:proofErr spellStart/spellEnd> around "Anzhelika" and
# "Rusyatinskaya", same as Word's own spellchecker would emit for the
# two words it doesn't recognise.
$p1 = $d.Paragraphs(1).Range
$p1.MoveEnd(1, -1)   # exclude the paragraph mark from the range

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/part.xml" pkg:contentType="application/xml">' +
       '<pkg:xmlData>' +
       '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:r><w:t xml:space="preserve">My name is </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>Anzhelika</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>Rusyatinskaya</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t>.</w:t></w:r>' +
       '</w:p>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$p1.InsertXML($xml)
